$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of row -> (DAMSLTag, DialogAct) in columns I and J
$updates = @{
    7   = @("sd", "Statement-non-opinion")
    51  = @("sv", "Statement-opinion")
    54  = @("sd", "Statement-non-opinion")
    56  = @("sd", "Statement-non-opinion")
    58  = @("b",  "Acknowledge (Backchannel)")
    59  = @("sd", "Statement-non-opinion")
    64  = @("sv", "Statement-opinion")
    68  = @("sv", "Statement-opinion")
    71  = @("ba", "Appreciation")
    77  = @("sd", "Statement-non-opinion")
    87  = @("ba", "Appreciation")
    92  = @("ba", "Appreciation")
    98  = @("%",  "Uninterpretable")
    131 = @("sd", "Statement-non-opinion")
    133 = @("sd", "Statement-non-opinion")
    144 = @("ba", "Appreciation")
    147 = @("sv", "Statement-opinion")
    149 = @("sd", "Statement-non-opinion")
    153 = @("sd", "Statement-non-opinion")
    159 = @("sv", "Statement-opinion")
    173 = @("ba", "Appreciation")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}
